# "Changed the base class and hooks"
# Update the Doctors worksheet: replace the 5 sample doctor rows with a new
# set of doctors/fields/experience/locations, and remove the now-unused
# "Status" column (F).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Doctors")

# New doctor rows: Name, Field, Experience, Practise Location
# (Surgeries List in column E is left untouched)
$rows = @(
    @{ Row=2; Name="Dr. Vikram Deshmukh";   Field="Dentist"; Experience="26 years experience overall"; Location="Nanded City,Pune" },
    @{ Row=3; Name="Dr. Gautami Phadke";    Field="Dentist"; Experience="26 years experience overall"; Location="Karve Nagar,Pune" },
    @{ Row=4; Name="Dr. Nisha R. Patil";    Field="Dentist"; Experience="21 years experience overall"; Location="Dhanori,Pune" },
    @{ Row=5; Name="Dr. Abhinav Misuriya";  Field="Dentist"; Experience="19 years experience overall"; Location="Hadapsar,Pune" },
    @{ Row=6; Name="Dr. Ritesh Khandelwal"; Field="Dentist"; Experience="25 years experience overall"; Location="Lonavala,Pune" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Name
    $ws.Cells.Item($r.Row, 2).Value = $r.Field
    $ws.Cells.Item($r.Row, 3).Value = $r.Experience
    $ws.Cells.Item($r.Row, 4).Value = $r.Location
}

# Remove the "Status" column header/data (column F) entirely.
$ws.Range("F1").ClearContents()
